$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Applied uniformly via a
# Text-format / Value / ClearFormats sequence so numeric-looking
# strings (e.g. "1.000", "306.24") are stored as text, matching the
# original inline-string cells, without leaving residual cell styles.
$cellUpdates = [ordered]@{
    "D2" = "23.469.69"
    "E2" = "  +1.25%  "
    "D3" = "1.639.45"
    "E3" = "  +2.37%  "
    "D4" = "1.000"
    "E4" = "  -0.10%  "
    "D6" = "306.24"
    "E6" = "  +0.98%  "
    "D7" = "0.3762"
    "E7" = "  -0.48%  "
    "D8" = "52.13"
    "E8" = "  +0.17%  "
    "D9" = "0.3641"
    "E9" = "  +0.76%  "
    "D10" = "1.263"
    "E10" = "  -0.33%  "
    "D11" = "0.08148"
    "E11" = "  +0.48%  "
    "D12" = "1.000"
    "E12" = "  -0.08%  "
    "D13" = "22.96"
    "E13" = "  +0.98%  "
    "D14" = "6.632"
    "E14" = "  +0.81%  "
    "D15" = "0.00001276"
    "E15" = "  +2.70%  "
    "D16" = "7.367"
    "E16" = "  -0.54%  "
    "D17" = "1.638.08"
    "E17" = "  +2.27%  "
    "D18" = "94.72"
    "E18" = "  +0.61%  "
    "D19" = "0.06904"
    "E19" = "  +0.36%  "
    "D20" = "18.19"
    "E20" = "  +0.64%  "
    "D21" = "6.545"
    "E21" = "  +0.01%  "
    "E22" = "  -0.02%  "
    "D23" = "23.482.44"
    "E23" = "  +1.30%  "
    "D24" = "12.78"
    "E24" = "  -1.43%  "
    "D25" = "3.084"
    "E25" = "  +3.44%  "
    "D26" = "2.420"
    "E26" = "  +0.86%  "
    "D27" = "21.26"
    "E27" = "  +0.13%  "
    "D28" = "150.83"
    "E28" = "  +0.92%  "
    "D29" = "5.348"
    "E29" = "  +2.26%  "
    "D30" = "136.88"
    "E30" = "  +2.15%  "
    "D31" = "2.311"
    "E31" = "  -3.19%  "
    "D32" = "1.819.63"
    "E32" = "  +2.22%  "
    "D33" = "6.783"
    "E33" = "  +0.45%  "
    "D34" = "0.9653"
    "E34" = "  -0.36%  "
    "D35" = "0.02832"
    "E35" = "  +4.35%  "
    "D36" = "10.32"
    "E36" = "  +0.38%  "
    "D37" = "0.07310"
    "E37" = "  -2.54%  "
    "D38" = "0.2530"
    "E38" = "  +1.18%  "
    "D39" = "0.08834"
    "E39" = "  +0.34%  "
    "D40" = "6.123"
    "E40" = "  +0.80%  "
    "E41" = "  +1.38%  "
    "D42" = "0.7101"
    "E42" = "  -0.10%  "
    "B43" = "Aptos"
    "C43" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D43" = "12.50"
    "E43" = "  +0.20%  "
    "B44" = "EnergySwap"
    "C44" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D44" = "16.24"
    "E44" = "  +3.94%  "
    "D45" = "0.6556"
    "E45" = "  +0.55%  "
    "D46" = "2.339"
    "E46" = "  +1.28%  "
    "D47" = "0.9997"
    "E47" = "  +0.02%  "
    "D48" = "4.016"
    "E48" = "  +0.03%  "
    "D49" = "0.07972"
    "E49" = "  +0.02%  "
    "D50" = "128.86"
    "E50" = "  -2.43%  "
    "D51" = "1.205"
    "E51" = "  +0.51%  "
}

foreach ($addr in $cellUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $cellUpdates[$addr]
    $ws.Range($addr).ClearFormats()
}

Write-Host "Applied $($cellUpdates.Count) cell updates to '$($ws.Name)'."
